$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
# Row 76
$ws.Cells.Item(76, 8).Value = 3876  # H76: 3628.6924 -> 3876
$ws.Cells.Item(76, 9).Value = 3417.889  # I76: 3348.25 -> 3417.889
$ws.Cells.Item(76, 10).Value = 4563.1665  # J76: 4077.4 -> 4563.1665
$ws.Cells.Item(76, 11).Value = 3417.889  # K76: 3348.25 -> 3417.889
$ws.Cells.Item(76, 12).Value = 4563.1665  # L76: 4077.4 -> 4563.1665
$ws.Cells.Item(76, 13).Value = -3102.889  # M76: -3033.25 -> -3102.889
$ws.Cells.Item(76, 14).Value = -5193.1665  # N76: -4707.4 -> -5193.1665
# Row 79
$ws.Cells.Item(79, 8).Value = 3876  # H79: 3628.6924 -> 3876
$ws.Cells.Item(79, 9).Value = 3417.889  # I79: 3348.25 -> 3417.889
$ws.Cells.Item(79, 10).Value = 4563.1665  # J79: 4077.4 -> 4563.1665
$ws.Cells.Item(79, 11).Value = 3417.889  # K79: 3348.25 -> 3417.889
$ws.Cells.Item(79, 12).Value = 4563.1665  # L79: 4077.4 -> 4563.1665
$ws.Cells.Item(79, 13).Value = -2325.889  # M79: -2256.25 -> -2325.889
$ws.Cells.Item(79, 14).Value = -6747.1665  # N79: -6261.4 -> -6747.1665
# Row 98
$ws.Cells.Item(98, 8).Value = 1294.9512  # H98: 1279.5476 -> 1294.9512
$ws.Cells.Item(98, 9).Value = 1264.725  # I98: 1264.725 -> 1264.725
$ws.Cells.Item(98, 10).Value = 2504  # J98: 1576 -> 2504
$ws.Cells.Item(98, 11).Value = 1264.725  # K98: 1264.725 -> 1264.725
$ws.Cells.Item(98, 12).Value = 2504  # L98: 1576 -> 2504
$ws.Cells.Item(98, 13).Value = 233.2750000000001  # M98: 233.2750000000001 -> 233.2750000000001
$ws.Cells.Item(98, 14).Value = -5500  # N98: -4572 -> -5500
# Row 101
$ws.Cells.Item(101, 8).Value = 438.46155  # H101: 438.6154 -> 438.46155
$ws.Cells.Item(101, 9).Value = 408.8  # I101: 409 -> 408.8
$ws.Cells.Item(101, 10).Value = 537.3333  # J101: 537.3333 -> 537.3333
$ws.Cells.Item(101, 11).Value = 1226.4  # K101: 1227 -> 1226.4
$ws.Cells.Item(101, 12).Value = 1611.9999  # L101: 1611.9999 -> 1611.9999
$ws.Cells.Item(101, 13).Value = 395.5999999999999  # M101: 395 -> 395.5999999999999
$ws.Cells.Item(101, 14).Value = -4855.9999  # N101: -4855.9999 -> -4855.9999
# Row 113
$ws.Cells.Item(113, 8).Value = 2350  # H113: 2287.25 -> 2350
$ws.Cells.Item(113, 9).Value = 2500  # I113: 2374.5 -> 2500
$ws.Cells.Item(113, 10).Value = 2200  # J113: 2200 -> 2200
$ws.Cells.Item(113, 11).Value = 2500  # K113: 2374.5 -> 2500
$ws.Cells.Item(113, 12).Value = 2200  # L113: 2200 -> 2200
$ws.Cells.Item(113, 13).Value = 754  # M113: 879.5 -> 754
$ws.Cells.Item(113, 14).Value = -8708  # N113: -8708 -> -8708
# Row 122
$ws.Cells.Item(122, 8).Value = 1294.9512  # H122: 1279.5476 -> 1294.9512
$ws.Cells.Item(122, 9).Value = 1264.725  # I122: 1264.725 -> 1264.725
$ws.Cells.Item(122, 10).Value = 2504  # J122: 1576 -> 2504
$ws.Cells.Item(122, 11).Value = 3794.175  # K122: 3794.175 -> 3794.175
$ws.Cells.Item(122, 12).Value = 7512  # L122: 4728 -> 7512
$ws.Cells.Item(122, 13).Value = -1344.175  # M122: -1344.175 -> -1344.175
$ws.Cells.Item(122, 14).Value = -12412  # N122: -9628 -> -12412
# Row 132
$ws.Cells.Item(132, 8).Value = 3641.0444  # H132: 3490.3618 -> 3641.0444
$ws.Cells.Item(132, 9).Value = 3626.475  # I132: 3458.5476 -> 3626.475
$ws.Cells.Item(132, 10).Value = 3757.6  # J132: 3757.6 -> 3757.6
$ws.Cells.Item(132, 11).Value = 10879.425  # K132: 10375.6428 -> 10879.425
$ws.Cells.Item(132, 12).Value = 11272.8  # L132: 11272.8 -> 11272.8
$ws.Cells.Item(132, 13).Value = -8349.424999999999  # M132: -7845.6428 -> -8349.424999999999
$ws.Cells.Item(132, 14).Value = -16332.8  # N132: -16332.8 -> -16332.8
# Row 138
$ws.Cells.Item(138, 8).Value = 3052.1133  # H138: 2918.224 -> 3052.1133
$ws.Cells.Item(138, 9).Value = 2884.1428  # I138: 2884.1428 -> 2884.1428
$ws.Cells.Item(138, 10).Value = 3112.4102  # J138: 2929.068 -> 3112.4102
$ws.Cells.Item(138, 11).Value = 8652.428400000001  # K138: 8652.428400000001 -> 8652.428400000001
$ws.Cells.Item(138, 12).Value = 9337.230599999999  # L138: 8787.204000000002 -> 9337.230599999999
$ws.Cells.Item(138, 13).Value = -3512.428400000001  # M138: -3512.428400000001 -> -3512.428400000001
$ws.Cells.Item(138, 14).Value = -19617.2306  # N138: -19067.204 -> -19617.2306

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Cells.Item(32, 8).Value = 9472329  # H32: 10419316 -> 9472329
$ws.Cells.Item(32, 9).Value = 5377863.5  # I32: 5748714.5 -> 5377863.5
$ws.Cells.Item(32, 10).Value = 19236056  # J32: 22732720 -> 19236056
$ws.Cells.Item(32, 11).Value = 5377863.5  # K32: 5748714.5 -> 5377863.5
$ws.Cells.Item(32, 12).Value = 19236056  # L32: 22732720 -> 19236056
$ws.Cells.Item(32, 13).Value = -5377576.5  # M32: -5748427.5 -> -5377576.5
$ws.Cells.Item(32, 14).Value = -19236630  # N32: -22733294 -> -19236630
# Row 45
$ws.Cells.Item(45, 8).Value = 4871.8184  # H45: 4528.6665 -> 4871.8184
$ws.Cells.Item(45, 9).Value = 3941.7144  # I45: 3543.25 -> 3941.7144
$ws.Cells.Item(45, 10).Value = 6499.5  # J45: 6499.5 -> 6499.5
$ws.Cells.Item(45, 11).Value = 3941.7144  # K45: 3543.25 -> 3941.7144
$ws.Cells.Item(45, 12).Value = 6499.5  # L45: 6499.5 -> 6499.5
$ws.Cells.Item(45, 13).Value = -3564.7144  # M45: -3166.25 -> -3564.7144
$ws.Cells.Item(45, 14).Value = -7253.5  # N45: -7253.5 -> -7253.5
# Row 74
$ws.Cells.Item(74, 8).Value = 1261.826  # H74: 1050.5667 -> 1261.826
$ws.Cells.Item(74, 9).Value = 1183.091  # I74: 983.5517 -> 1183.091
$ws.Cells.Item(74, 10).Value = 2994  # J74: 2994 -> 2994
$ws.Cells.Item(74, 11).Value = 1183.091  # K74: 983.5517 -> 1183.091
$ws.Cells.Item(74, 12).Value = 2994  # L74: 2994 -> 2994
$ws.Cells.Item(74, 13).Value = -309.0909999999999  # M74: -109.5517 -> -309.0909999999999
$ws.Cells.Item(74, 14).Value = -4742  # N74: -4742 -> -4742
# Row 77
$ws.Cells.Item(77, 8).Value = 1261.826  # H77: 1050.5667 -> 1261.826
$ws.Cells.Item(77, 9).Value = 1183.091  # I77: 983.5517 -> 1183.091
$ws.Cells.Item(77, 10).Value = 2994  # J77: 2994 -> 2994
$ws.Cells.Item(77, 11).Value = 5915.455  # K77: 4917.7585 -> 5915.455
$ws.Cells.Item(77, 12).Value = 14970  # L77: 14970 -> 14970
$ws.Cells.Item(77, 13).Value = -1547.455  # M77: -549.7584999999999 -> -1547.455
$ws.Cells.Item(77, 14).Value = -23706  # N77: -23706 -> -23706
# Row 88
$ws.Cells.Item(88, 8).Value = 1137.6666  # H88: 1206.5 -> 1137.6666
$ws.Cells.Item(88, 9).Value = 1203  # I88: 1406 -> 1203
$ws.Cells.Item(88, 10).Value = 1007  # J88: 1007 -> 1007
$ws.Cells.Item(88, 11).Value = 1203  # K88: 1406 -> 1203
$ws.Cells.Item(88, 12).Value = 1007  # L88: 1007 -> 1007
$ws.Cells.Item(88, 13).Value = -797  # M88: -1000 -> -797
$ws.Cells.Item(88, 14).Value = -1819  # N88: -1819 -> -1819
# Row 91
$ws.Cells.Item(91, 8).Value = 1137.6666  # H91: 1206.5 -> 1137.6666
$ws.Cells.Item(91, 9).Value = 1203  # I91: 1406 -> 1203
$ws.Cells.Item(91, 10).Value = 1007  # J91: 1007 -> 1007
$ws.Cells.Item(91, 11).Value = 1203  # K91: 1406 -> 1203
$ws.Cells.Item(91, 12).Value = 1007  # L91: 1007 -> 1007
$ws.Cells.Item(91, 13).Value = 201  # M91: -2 -> 201
$ws.Cells.Item(91, 14).Value = -3815  # N91: -3815 -> -3815
# Row 97
$ws.Cells.Item(97, 8).Value = 1197.5312  # H97: 1230.0322 -> 1197.5312
$ws.Cells.Item(97, 9).Value = 1009.2857  # I97: 1050.25 -> 1009.2857
$ws.Cells.Item(97, 10).Value = 1556.909  # J97: 1556.909 -> 1556.909
$ws.Cells.Item(97, 11).Value = 1009.2857  # K97: 1050.25 -> 1009.2857
$ws.Cells.Item(97, 12).Value = 1556.909  # L97: 1556.909 -> 1556.909
$ws.Cells.Item(97, 13).Value = -513.2857  # M97: -554.25 -> -513.2857
$ws.Cells.Item(97, 14).Value = -2548.909  # N97: -2548.909 -> -2548.909
# Row 134
$ws.Cells.Item(134, 8).Value = 0  # H134: 55000 -> 0
$ws.Cells.Item(134, 9).Value = 0  # I134: 0 -> 0
$ws.Cells.Item(134, 10).Value = 0  # J134: 55000 -> 0
$ws.Cells.Item(134, 11).Value = 0  # K134: 0 -> 0
$ws.Cells.Item(134, 12).ClearContents()  # L134: 55000 -> (removed)
$ws.Cells.Item(134, 14).Value = 0  # N134: -65140 -> 0

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
# Row 61
$ws.Cells.Item(61, 8).Value = 30000  # H61: 0 -> 30000
$ws.Cells.Item(61, 9).Value = 0  # I61: 0 -> 0
$ws.Cells.Item(61, 10).Value = 30000  # J61: 0 -> 30000
$ws.Cells.Item(61, 11).Value = 0  # K61: 0 -> 0
$ws.Cells.Item(61, 12).Value = 30000  # L61: 0 -> 30000
$ws.Cells.Item(61, 14).Value = -30626  # N61: None -> -30626
# Row 137
$ws.Cells.Item(137, 8).Value = 44999.5  # H137: 0 -> 44999.5
$ws.Cells.Item(137, 9).Value = 0  # I137: 0 -> 0
$ws.Cells.Item(137, 10).Value = 44999.5  # J137: 0 -> 44999.5
$ws.Cells.Item(137, 11).Value = 0  # K137: 0 -> 0
$ws.Cells.Item(137, 12).Value = 44999.5  # L137: 0 -> 44999.5
$ws.Cells.Item(137, 14).Value = -55199.5  # N137: None -> -55199.5
# Row 138
$ws.Cells.Item(138, 8).Value = 80000  # H138: 100000 -> 80000
$ws.Cells.Item(138, 9).Value = 80000  # I138: 80000 -> 80000
$ws.Cells.Item(138, 10).Value = 0  # J138: 120000 -> 0
$ws.Cells.Item(138, 11).Value = 80000  # K138: 80000 -> 80000
$ws.Cells.Item(138, 12).Value = 0  # L138: 120000 -> 0
$ws.Cells.Item(138, 13).Value = -74860  # M138: -74860 -> -74860
$ws.Cells.Item(138, 14).ClearContents()  # N138: -130280 -> (removed)
# Row 139
$ws.Cells.Item(139, 8).Value = 120000.336  # H139: 112000.2 -> 120000.336
$ws.Cells.Item(139, 9).Value = 0  # I139: 0 -> 0
$ws.Cells.Item(139, 10).Value = 120000.336  # J139: 112000.2 -> 120000.336
$ws.Cells.Item(139, 11).Value = 0  # K139: 0 -> 0
$ws.Cells.Item(139, 12).Value = 120000.336  # L139: 112000.2 -> 120000.336
$ws.Cells.Item(139, 14).Value = -130280.336  # N139: -122280.2 -> -130280.336

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
# Row 31
$ws.Cells.Item(31, 8).Value = 3359.25  # H31: 4003 -> 3359.25
$ws.Cells.Item(31, 9).Value = 1862.9231  # I31: 2342.6667 -> 1862.9231
$ws.Cells.Item(31, 10).Value = 6138.143  # J31: 5870.875 -> 6138.143
$ws.Cells.Item(31, 11).Value = 1862.9231  # K31: 2342.6667 -> 1862.9231
$ws.Cells.Item(31, 12).Value = 6138.143  # L31: 5870.875 -> 6138.143
$ws.Cells.Item(31, 13).Value = -1567.9231  # M31: -2047.6667 -> -1567.9231
$ws.Cells.Item(31, 14).Value = -6728.143  # N31: -6460.875 -> -6728.143
# Row 34
$ws.Cells.Item(34, 8).Value = 3359.25  # H34: 4003 -> 3359.25
$ws.Cells.Item(34, 9).Value = 1862.9231  # I34: 2342.6667 -> 1862.9231
$ws.Cells.Item(34, 10).Value = 6138.143  # J34: 5870.875 -> 6138.143
$ws.Cells.Item(34, 11).Value = 1862.9231  # K34: 2342.6667 -> 1862.9231
$ws.Cells.Item(34, 12).Value = 6138.143  # L34: 5870.875 -> 6138.143
$ws.Cells.Item(34, 13).Value = -1660.9231  # M34: -2140.6667 -> -1660.9231
$ws.Cells.Item(34, 14).Value = -6274.875  # N34: -6274.875 -> -6274.875
# Row 131
$ws.Cells.Item(131, 8).Value = 163000  # H131: 162994.5 -> 163000
$ws.Cells.Item(131, 9).Value = 0  # I131: 0 -> 0
$ws.Cells.Item(131, 10).Value = 163000  # J131: 162994.5 -> 163000
$ws.Cells.Item(131, 11).Value = 0  # K131: 0 -> 0
$ws.Cells.Item(131, 12).Value = 163000  # L131: 162994.5 -> 163000
$ws.Cells.Item(131, 14).Value = -173080  # N131: -173074.5 -> -173080
# Row 134
$ws.Cells.Item(134, 8).Value = 2806.9092  # H134: 2748.8696 -> 2806.9092
$ws.Cells.Item(134, 9).Value = 2057.8667  # I134: 2021.25 -> 2057.8667
$ws.Cells.Item(134, 10).Value = 4412  # J134: 4412 -> 4412
$ws.Cells.Item(134, 11).Value = 6173.6001  # K134: 6063.75 -> 6173.6001
$ws.Cells.Item(134, 12).Value = 13236  # L134: 13236 -> 13236
$ws.Cells.Item(134, 13).Value = -3638.6001  # M134: -3528.75 -> -3638.6001
$ws.Cells.Item(134, 14).Value = -18306  # N134: -18306 -> -18306

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
# Row 5
$ws.Cells.Item(5, 8).Value = 1163.8572  # H5: 1115.625 -> 1163.8572
$ws.Cells.Item(5, 9).Value = 599  # I5: 598.5 -> 599
$ws.Cells.Item(5, 10).Value = 1258  # J5: 1288 -> 1258
$ws.Cells.Item(5, 11).Value = 1797  # K5: 1795.5 -> 1797
$ws.Cells.Item(5, 12).Value = 3774  # L5: 3864 -> 3774
$ws.Cells.Item(5, 13).Value = -1685  # M5: -1683.5 -> -1685
$ws.Cells.Item(5, 14).Value = -3998  # N5: -4088 -> -3998
# Row 75
$ws.Cells.Item(75, 8).Value = 0  # H75: 1542.75 -> 0
$ws.Cells.Item(75, 9).Value = 0  # I75: 398.25 -> 0
$ws.Cells.Item(75, 10).Value = 0  # J75: 2687.25 -> 0
$ws.Cells.Item(75, 11).Value = 0  # K75: 1194.75 -> 0
$ws.Cells.Item(75, 12).ClearContents()  # L75: 8061.75 -> (removed)
$ws.Cells.Item(75, 13).ClearContents()  # M75: -196.75 -> (removed)
$ws.Cells.Item(75, 14).Value = 0  # N75: -10057.75 -> 0
# Row 78
$ws.Cells.Item(78, 8).Value = 0  # H78: 1542.75 -> 0
$ws.Cells.Item(78, 9).Value = 0  # I78: 398.25 -> 0
$ws.Cells.Item(78, 10).Value = 0  # J78: 2687.25 -> 0
$ws.Cells.Item(78, 11).Value = 0  # K78: 3584.25 -> 0
$ws.Cells.Item(78, 12).ClearContents()  # L78: 24185.25 -> (removed)
$ws.Cells.Item(78, 13).ClearContents()  # M78: 1407.75 -> (removed)
$ws.Cells.Item(78, 14).Value = 0  # N78: -34169.25 -> 0
# Row 82
$ws.Cells.Item(82, 8).Value = 1499  # H82: 0 -> 1499
$ws.Cells.Item(82, 9).Value = 1499  # I82: 0 -> 1499
$ws.Cells.Item(82, 10).Value = 0  # J82: 0 -> 0
$ws.Cells.Item(82, 11).Value = 4497  # K82: 0 -> 4497
$ws.Cells.Item(82, 12).Value = 0  # L82: 0 -> 0
$ws.Cells.Item(82, 13).Value = -4091  # M82: None -> -4091
# Row 85
$ws.Cells.Item(85, 8).Value = 1499  # H85: 0 -> 1499
$ws.Cells.Item(85, 9).Value = 1499  # I85: 0 -> 1499
$ws.Cells.Item(85, 10).Value = 0  # J85: 0 -> 0
$ws.Cells.Item(85, 11).Value = 4497  # K85: 0 -> 4497
$ws.Cells.Item(85, 12).Value = 0  # L85: 0 -> 0
$ws.Cells.Item(85, 13).Value = -3093  # M85: None -> -3093
# Row 107
$ws.Cells.Item(107, 8).Value = 606.05  # H107: 584.7727 -> 606.05
$ws.Cells.Item(107, 9).Value = 1166.3334  # I107: 1166.3334 -> 1166.3334
$ws.Cells.Item(107, 10).Value = 507.17648  # J107: 492.94736 -> 507.17648
$ws.Cells.Item(107, 11).Value = 3499.0002  # K107: 3499.0002 -> 3499.0002
$ws.Cells.Item(107, 12).Value = 1521.52944  # L107: 1478.84208 -> 1521.52944
$ws.Cells.Item(107, 13).Value = -1579.0002  # M107: -1579.0002 -> -1579.0002
$ws.Cells.Item(107, 14).Value = -5361.52944  # N107: -5318.84208 -> -5361.52944
# Row 121
$ws.Cells.Item(121, 8).Value = 12626142  # H121: 27286230 -> 12626142
$ws.Cells.Item(121, 9).Value = 16667355  # I121: 20000722 -> 16667355
$ws.Cells.Item(121, 10).Value = 502501.5  # J121: 45500000 -> 502501.5
$ws.Cells.Item(121, 11).Value = 50002065  # K121: 60002166 -> 50002065
$ws.Cells.Item(121, 12).Value = 1507504.5  # L121: 136500000 -> 1507504.5
$ws.Cells.Item(121, 13).Value = -50000755  # M121: -60000856 -> -50000755
$ws.Cells.Item(121, 14).Value = -1510124.5  # N121: -136502620 -> -1510124.5
# Row 128
$ws.Cells.Item(128, 8).Value = 91644.39999999999  # H128: 88370.164 -> 91644.39999999999
$ws.Cells.Item(128, 9).Value = 91644.39999999999  # I128: 88370.164 -> 91644.39999999999
$ws.Cells.Item(128, 10).Value = 0  # J128: 0 -> 0
$ws.Cells.Item(128, 11).Value = 274933.2  # K128: 265110.492 -> 274933.2
$ws.Cells.Item(128, 12).Value = 0  # L128: 0 -> 0
$ws.Cells.Item(128, 13).Value = -269953.2  # M128: -260130.492 -> -269953.2
# Row 131
$ws.Cells.Item(131, 8).Value = 1694.4736  # H131: 1732 -> 1694.4736
$ws.Cells.Item(131, 9).Value = 1465.625  # I131: 1585.8334 -> 1465.625
$ws.Cells.Item(131, 10).Value = 1755.5  # J131: 1764.4814 -> 1755.5
$ws.Cells.Item(131, 11).Value = 4396.875  # K131: 4757.5002 -> 4396.875
$ws.Cells.Item(131, 12).Value = 5266.5  # L131: 5293.4442 -> 5266.5
$ws.Cells.Item(131, 13).Value = 643.125  # M131: 282.4997999999996 -> 643.125
$ws.Cells.Item(131, 14).Value = -15346.5  # N131: -15373.4442 -> -15346.5
# Row 132
$ws.Cells.Item(132, 8).Value = 1557.7142  # H132: 1503.7826 -> 1557.7142
$ws.Cells.Item(132, 9).Value = 979.5  # I132: 971.1 -> 979.5
$ws.Cells.Item(132, 10).Value = 1913.5385  # J132: 1913.5385 -> 1913.5385
$ws.Cells.Item(132, 11).Value = 8815.5  # K132: 8739.9 -> 8815.5
$ws.Cells.Item(132, 12).Value = 17221.8465  # L132: 17221.8465 -> 17221.8465
$ws.Cells.Item(132, 13).Value = -6285.5  # M132: -6209.9 -> -6285.5
$ws.Cells.Item(132, 14).Value = -22281.8465  # N132: -22281.8465 -> -22281.8465
# Row 135
$ws.Cells.Item(135, 8).Value = 1163.8572  # H135: 1115.625 -> 1163.8572
$ws.Cells.Item(135, 9).Value = 599  # I135: 598.5 -> 599
$ws.Cells.Item(135, 10).Value = 1258  # J135: 1288 -> 1258
$ws.Cells.Item(135, 11).Value = 5391  # K135: 5386.5 -> 5391
$ws.Cells.Item(135, 12).Value = 11322  # L135: 11592 -> 11322
$ws.Cells.Item(135, 13).Value = -2856  # M135: -2851.5 -> -2856
$ws.Cells.Item(135, 14).Value = -16392  # N135: -16662 -> -16392

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
# Row 132
$ws.Cells.Item(132, 8).Value = 4813.857  # H132: 2654.0435 -> 4813.857
$ws.Cells.Item(132, 9).Value = 3864.7778  # I132: 1778.375 -> 3864.7778
$ws.Cells.Item(132, 10).Value = 6522.2  # J132: 4655.5713 -> 6522.2
$ws.Cells.Item(132, 11).Value = 11594.3334  # K132: 5335.125 -> 11594.3334
$ws.Cells.Item(132, 12).Value = 19566.6  # L132: 13966.7139 -> 19566.6
$ws.Cells.Item(132, 13).Value = -9064.3334  # M132: -2805.125 -> -9064.3334
$ws.Cells.Item(132, 14).Value = -24626.6  # N132: -19026.7139 -> -24626.6

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
# Row 7
$ws.Cells.Item(7, 8).Value = 4041.7368  # H7: 4120.684 -> 4041.7368
$ws.Cells.Item(7, 9).Value = 3191.4546  # I7: 3327.818 -> 3191.4546
$ws.Cells.Item(7, 10).Value = 5210.875  # J7: 5210.875 -> 5210.875
$ws.Cells.Item(7, 11).Value = 3191.4546  # K7: 3327.818 -> 3191.4546
$ws.Cells.Item(7, 12).Value = 5210.875  # L7: 5210.875 -> 5210.875
$ws.Cells.Item(7, 13).Value = -3079.4546  # M7: -3215.818 -> -3079.4546
$ws.Cells.Item(7, 14).Value = -5434.875  # N7: -5434.875 -> -5434.875
# Row 16
$ws.Cells.Item(16, 8).Value = 1567.3572  # H16: 1746.0834 -> 1567.3572
$ws.Cells.Item(16, 9).Value = 1534.0769  # I16: 1723 -> 1534.0769
$ws.Cells.Item(16, 10).Value = 2000  # J16: 2000 -> 2000
$ws.Cells.Item(16, 11).Value = 1534.0769  # K16: 1723 -> 1534.0769
$ws.Cells.Item(16, 12).Value = 2000  # L16: 2000 -> 2000
$ws.Cells.Item(16, 13).Value = -1364.0769  # M16: -1553 -> -1364.0769
$ws.Cells.Item(16, 14).Value = -2340  # N16: -2340 -> -2340
# Row 62
$ws.Cells.Item(62, 8).Value = 25000  # H62: 0 -> 25000
$ws.Cells.Item(62, 9).Value = 25000  # I62: 0 -> 25000
$ws.Cells.Item(62, 10).Value = 0  # J62: 0 -> 0
$ws.Cells.Item(62, 11).Value = 25000  # K62: 0 -> 25000
$ws.Cells.Item(62, 12).Value = 0  # L62: 0 -> 0
$ws.Cells.Item(62, 13).Value = -24376  # M62: None -> -24376
# Row 65
$ws.Cells.Item(65, 8).Value = 25000  # H65: 0 -> 25000
$ws.Cells.Item(65, 9).Value = 25000  # I65: 0 -> 25000
$ws.Cells.Item(65, 10).Value = 0  # J65: 0 -> 0
$ws.Cells.Item(65, 11).Value = 75000  # K65: 0 -> 75000
$ws.Cells.Item(65, 12).Value = 0  # L65: 0 -> 0
$ws.Cells.Item(65, 13).Value = -71880  # M65: None -> -71880
# Row 82
$ws.Cells.Item(82, 8).Value = 3230.7778  # H82: 3375.2354 -> 3230.7778
$ws.Cells.Item(82, 9).Value = 2559.5  # I82: 2916.4 -> 2559.5
$ws.Cells.Item(82, 10).Value = 3566.4167  # J82: 3566.4167 -> 3566.4167
$ws.Cells.Item(82, 11).Value = 2559.5  # K82: 2916.4 -> 2559.5
$ws.Cells.Item(82, 12).Value = 3566.4167  # L82: 3566.4167 -> 3566.4167
$ws.Cells.Item(82, 13).Value = -2198.5  # M82: -2555.4 -> -2198.5
$ws.Cells.Item(82, 14).Value = -4288.4167  # N82: -4288.4167 -> -4288.4167
# Row 85
$ws.Cells.Item(85, 8).Value = 3230.7778  # H85: 3375.2354 -> 3230.7778
$ws.Cells.Item(85, 9).Value = 2559.5  # I85: 2916.4 -> 2559.5
$ws.Cells.Item(85, 10).Value = 3566.4167  # J85: 3566.4167 -> 3566.4167
$ws.Cells.Item(85, 11).Value = 2559.5  # K85: 2916.4 -> 2559.5
$ws.Cells.Item(85, 12).Value = 3566.4167  # L85: 3566.4167 -> 3566.4167
$ws.Cells.Item(85, 13).Value = -1311.5  # M85: -1668.4 -> -1311.5
$ws.Cells.Item(85, 14).Value = -6062.4167  # N85: -6062.4167 -> -6062.4167
# Row 100
$ws.Cells.Item(100, 8).Value = 1835.8334  # H100: 1809.9231 -> 1835.8334
$ws.Cells.Item(100, 9).Value = 1888.625  # I100: 1845.3334 -> 1888.625
$ws.Cells.Item(100, 10).Value = 1730.25  # J100: 1730.25 -> 1730.25
$ws.Cells.Item(100, 11).Value = 1888.625  # K100: 1845.3334 -> 1888.625
$ws.Cells.Item(100, 12).Value = 1730.25  # L100: 1730.25 -> 1730.25
$ws.Cells.Item(100, 13).Value = -1347.625  # M100: -1304.3334 -> -1347.625
$ws.Cells.Item(100, 14).Value = -2812.25  # N100: -2812.25 -> -2812.25
# Row 126
$ws.Cells.Item(126, 8).Value = 4041.7368  # H126: 4120.684 -> 4041.7368
$ws.Cells.Item(126, 9).Value = 3191.4546  # I126: 3327.818 -> 3191.4546
$ws.Cells.Item(126, 10).Value = 5210.875  # J126: 5210.875 -> 5210.875
$ws.Cells.Item(126, 11).Value = 9574.363799999999  # K126: 9983.454000000002 -> 9574.363799999999
$ws.Cells.Item(126, 12).Value = 15632.625  # L126: 15632.625 -> 15632.625
$ws.Cells.Item(126, 13).Value = -7104.363799999999  # M126: -7513.454000000002 -> -7104.363799999999
$ws.Cells.Item(126, 14).Value = -20572.625  # N126: -20572.625 -> -20572.625

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
# Row 132
$ws.Cells.Item(132, 8).Value = 2065  # H132: 2018.2778 -> 2065
$ws.Cells.Item(132, 9).Value = 1557.1957  # I132: 1532.0851 -> 1557.1957
$ws.Cells.Item(132, 10).Value = 4984.875  # J132: 5282.7144 -> 4984.875
$ws.Cells.Item(132, 11).Value = 4671.5871  # K132: 4596.2553 -> 4671.5871
$ws.Cells.Item(132, 12).Value = 14954.625  # L132: 15848.1432 -> 14954.625
$ws.Cells.Item(132, 13).Value = -2141.5871  # M132: -2066.2553 -> -2141.5871
$ws.Cells.Item(132, 14).Value = -20014.625  # N132: -20908.1432 -> -20014.625
# Row 136
$ws.Cells.Item(136, 8).Value = 2443.0557  # H136: 2560.585 -> 2443.0557
$ws.Cells.Item(136, 9).Value = 1517.9231  # I136: 1657.5 -> 1517.9231
$ws.Cells.Item(136, 10).Value = 4848.4  # J136: 4848.4 -> 4848.4
$ws.Cells.Item(136, 11).Value = 4553.7693  # K136: 4972.5 -> 4553.7693
$ws.Cells.Item(136, 12).Value = 14545.2  # L136: 14545.2 -> 14545.2
$ws.Cells.Item(136, 13).Value = -2003.7693  # M136: -2422.5 -> -2003.7693
$ws.Cells.Item(136, 14).Value = -19645.2  # N136: -19645.2 -> -19645.2
